$d = $word.ActiveDocument

# --- Change 1: merge the "Sexta Reunión (12/04/2022) " heading's runs into a single run ---
# (the text itself is unchanged; re-running Find/Replace over the whole phrase collapses
# the four pre-existing runs into one run, matching the target XML)
$d.Content.Find.Execute("Sexta Reunión (12/04/2022) ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Sexta Reunión (12/04/2022) ", 2) | Out-Null

# --- Change 2: append the new "Séptima Reunión" block after the last paragraph ---
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = @"
<w:p $ns><w:pPr><w:jc w:val="center"/></w:pPr></w:p>
<w:p $ns><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>Séptima</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Reunión (1</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>9</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">/04/2022) </w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Debate sobre cuestiones referentes a</w:t></w:r><w:r><w:t xml:space="preserve"> la presentación </w:t></w:r><w:r><w:t>del diseño del videojuego.</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">Puesta en común de todos los avances desarrollados a lo largo de la </w:t></w:r><w:r><w:t>fase 2</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Exposición de buenas practicas de desarrollo y mantenimiento de contenido y archivos.</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">Inicio de la fase </w:t></w:r><w:r><w:t>3</w:t></w:r><w:r><w:t xml:space="preserve"> de la planificación de desarrollo del proyecto.</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="center"/></w:pPr></w:p>
"@

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($newXml)
